# Fruta / hortaliza, semanal
# Rows 3-10 get their weekly observation data (date, volume, prices,
# unit, origin, $/Kg, Kg/unidad) reshuffled across the block. The
# descriptive columns (A,B,C,E,F,G,H,I,J,K,L) are identical for every
# row in the block, so only D and M:T actually change per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to write into rows 3-10, column by column:
# D = Fecha (date serial), M = Volumen, N = Precio minimo,
# O = Precio maximo, P = Precio promedio ponderado,
# Q = Unidad de comercializacion, R = Origen, S = Precio $/Kg,
# T = Kg / unidad
$rows = @{
    3  = @{ D = 44187; M = 80;  N = 2800; O = 3000; P = 2900; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";    S = 1450; T = 2 }
    4  = @{ D = 44187; M = 65;  N = 1400; O = 1500; P = 1446; Q = "`$/envase 1 kilo";    R = "Provincia de Diguillín"; S = 1446; T = 1 }
    5  = @{ D = 44539; M = 200; N = 3800; O = 4000; P = 3900; Q = "`$/bandeja 2 kilos"; R = "Región del Maule";       S = 1950; T = 2 }
    6  = @{ D = 44594; M = 120; N = 2500; O = 2800; P = 2650; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";    S = 1325; T = 2 }
    7  = @{ D = 44174; M = 150; N = 3700; O = 3800; P = 3747; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";    S = 1874; T = 2 }
    8  = @{ D = 44596; M = 120; N = 2500; O = 2700; P = 2600; Q = "`$/bandeja 2 kilos"; R = "Provincia de Linares";    S = 1300; T = 2 }
    9  = @{ D = 44181; M = 65;  N = 3600; O = 3800; P = 3692; Q = "`$/bandeja 2 kilos"; R = "Provincia de Diguillín"; S = 1846; T = 2 }
    10 = @{ D = 44181; M = 80;  N = 1800; O = 2000; P = 1875; Q = "`$/envase 1 kilo";    R = "Provincia de Diguillín"; S = 1875; T = 1 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
